$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after existing row 901 (shifts old rows
# 902-934 down to 904-936, carrying formatting/styles down with them).
$ws.Rows("902:903").Insert()

# Row 902 (new): Fruta, Terminal Hortofrutícola Agro Chillán - Plátano, Pintón
$ws.Cells.Item(902, 1).Value  = 7
$ws.Cells.Item(902, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(902, 3).Value  = "Ñuble"
$ws.Cells.Item(902, 4).Value  = 45075
$ws.Cells.Item(902, 5).Value  = 16
$ws.Cells.Item(902, 6).Value  = "Fruta"
$ws.Cells.Item(902, 7).Value  = 100108
$ws.Cells.Item(902, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(902, 9).Value  = 100108006
$ws.Cells.Item(902, 10).Value = "Plátano"
$ws.Cells.Item(902, 11).Value = "Sin especificar"
$ws.Cells.Item(902, 12).Value = "Pintón"
$ws.Cells.Item(902, 13).Value = 150
$ws.Cells.Item(902, 14).Value = 15000
$ws.Cells.Item(902, 15).Value = 15000
$ws.Cells.Item(902, 16).Value = 15000
$ws.Cells.Item(902, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(902, 18).Value = "Ecuador"
$ws.Cells.Item(902, 19).Value = 750
$ws.Cells.Item(902, 20).Value = 20

# Row 903 (new): Fruta, Terminal Hortofrutícola Agro Chillán - Plátano, Primera Pintón
$ws.Cells.Item(903, 1).Value  = 7
$ws.Cells.Item(903, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(903, 3).Value  = "Ñuble"
$ws.Cells.Item(903, 4).Value  = 45075
$ws.Cells.Item(903, 5).Value  = 16
$ws.Cells.Item(903, 6).Value  = "Fruta"
$ws.Cells.Item(903, 7).Value  = 100108
$ws.Cells.Item(903, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(903, 9).Value  = 100108006
$ws.Cells.Item(903, 10).Value = "Plátano"
$ws.Cells.Item(903, 11).Value = "Sin especificar"
$ws.Cells.Item(903, 12).Value = "Primera Pintón"
$ws.Cells.Item(903, 13).Value = 100
$ws.Cells.Item(903, 14).Value = 16000
$ws.Cells.Item(903, 15).Value = 16000
$ws.Cells.Item(903, 16).Value = 16000
$ws.Cells.Item(903, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(903, 18).Value = "Ecuador"
$ws.Cells.Item(903, 19).Value = 800
$ws.Cells.Item(903, 20).Value = 20
